$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '23.005.34'; Text = $true }
    @{ Cell = 'E2'; Value = '  -3.74%  '; Text = $false }
    @{ Cell = 'D3'; Value = '1.602.27'; Text = $true }
    @{ Cell = 'E3'; Value = '  -2.85%  '; Text = $false }
    @{ Cell = 'D4'; Value = '1.002'; Text = $true }
    @{ Cell = 'E4'; Value = '  +0.12%  '; Text = $false }
    @{ Cell = 'E5'; Value = '  +0.09%  '; Text = $false }
    @{ Cell = 'D6'; Value = '301.09'; Text = $true }
    @{ Cell = 'E6'; Value = '  -3.08%  '; Text = $false }
    @{ Cell = 'D7'; Value = '0.3776'; Text = $true }
    @{ Cell = 'E7'; Value = '  -3.04%  '; Text = $false }
    @{ Cell = 'D8'; Value = '0.3623'; Text = $true }
    @{ Cell = 'E8'; Value = '  -5.61%  '; Text = $false }
    @{ Cell = 'D9'; Value = '49.57'; Text = $true }
    @{ Cell = 'E9'; Value = '  -2.80%  '; Text = $false }
    @{ Cell = 'E10'; Value = '  -6.25%  '; Text = $false }
    @{ Cell = 'E11'; Value = '  +0.15%  '; Text = $false }
    @{ Cell = 'D12'; Value = '0.08120'; Text = $true }
    @{ Cell = 'E12'; Value = '  -3.85%  '; Text = $false }
    @{ Cell = 'E13'; Value = '  -4.22%  '; Text = $false }
    @{ Cell = 'D14'; Value = '6.588'; Text = $true }
    @{ Cell = 'E14'; Value = '  -6.19%  '; Text = $false }
    @{ Cell = 'D15'; Value = '7.388'; Text = $true }
    @{ Cell = 'E15'; Value = '  -6.70%  '; Text = $false }
    @{ Cell = 'D16'; Value = '0.00001241'; Text = $true }
    @{ Cell = 'E16'; Value = '  -5.68%  '; Text = $false }
    @{ Cell = 'D17'; Value = '1.598.48'; Text = $true }
    @{ Cell = 'E17'; Value = '  -3.07%  '; Text = $false }
    @{ Cell = 'D18'; Value = '92.08'; Text = $true }
    @{ Cell = 'E18'; Value = '  -1.96%  '; Text = $false }
    @{ Cell = 'D19'; Value = '0.06871'; Text = $true }
    @{ Cell = 'E19'; Value = '  -1.35%  '; Text = $false }
    @{ Cell = 'E20'; Value = '  -6.56%  '; Text = $false }
    @{ Cell = 'D21'; Value = '6.561'; Text = $true }
    @{ Cell = 'E21'; Value = '  -5.43%  '; Text = $false }
    @{ Cell = 'D22'; Value = '1.001'; Text = $true }
    @{ Cell = 'E22'; Value = '  +0.08%  '; Text = $false }
    @{ Cell = 'D23'; Value = '13.16'; Text = $true }
    @{ Cell = 'E23'; Value = '  -3.44%  '; Text = $false }
    @{ Cell = 'D24'; Value = '23.013.87'; Text = $true }
    @{ Cell = 'E24'; Value = '  -3.72%  '; Text = $false }
    @{ Cell = 'D25'; Value = '2.366'; Text = $true }
    @{ Cell = 'E25'; Value = '  -3.07%  '; Text = $false }
    @{ Cell = 'D26'; Value = '2.803'; Text = $true }
    @{ Cell = 'E26'; Value = '  -3.52%  '; Text = $false }
    @{ Cell = 'D27'; Value = '21.06'; Text = $true }
    @{ Cell = 'E27'; Value = '  -4.12%  '; Text = $false }
    @{ Cell = 'D28'; Value = '150.48'; Text = $true }
    @{ Cell = 'E28'; Value = '  -2.41%  '; Text = $false }
    @{ Cell = 'D29'; Value = '5.252'; Text = $true }
    @{ Cell = 'E29'; Value = '  -2.61%  '; Text = $false }
    @{ Cell = 'D30'; Value = '133.46'; Text = $true }
    @{ Cell = 'E30'; Value = '  -2.73%  '; Text = $false }
    @{ Cell = 'D31'; Value = '2.312'; Text = $true }
    @{ Cell = 'E31'; Value = '  -6.84%  '; Text = $false }
    @{ Cell = 'D32'; Value = '6.806'; Text = $true }
    @{ Cell = 'E32'; Value = '  -11.80%  '; Text = $false }
    @{ Cell = 'D33'; Value = '1.778.47'; Text = $true }
    @{ Cell = 'E33'; Value = '  -2.84%  '; Text = $false }
    @{ Cell = 'D34'; Value = '0.9609'; Text = $true }
    @{ Cell = 'E34'; Value = '  -2.88%  '; Text = $false }
    @{ Cell = 'D35'; Value = '0.07627'; Text = $true }
    @{ Cell = 'E35'; Value = '  -6.05%  '; Text = $false }
    @{ Cell = 'D36'; Value = '10.38'; Text = $true }
    @{ Cell = 'E36'; Value = '  -0.71%  '; Text = $false }
    @{ Cell = 'D37'; Value = '6.300'; Text = $true }
    @{ Cell = 'E37'; Value = '  -5.75%  '; Text = $false }
    @{ Cell = 'B38'; Value = 'Algorand'; Text = $false }
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; Text = $false }
    @{ Cell = 'D38'; Value = '0.2537'; Text = $true }
    @{ Cell = 'E38'; Value = '  -5.35%  '; Text = $false }
    @{ Cell = 'B39'; Value = 'VeChain'; Text = $false }
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; Text = $false }
    @{ Cell = 'D39'; Value = '0.02701'; Text = $true }
    @{ Cell = 'E39'; Value = '  -7.53%  '; Text = $false }
    @{ Cell = 'D40'; Value = '0.08855'; Text = $true }
    @{ Cell = 'E40'; Value = '  -2.89%  '; Text = $false }
    @{ Cell = 'E41'; Value = '  -4.03%  '; Text = $false }
    @{ Cell = 'D42'; Value = '0.7051'; Text = $true }
    @{ Cell = 'E42'; Value = '  -6.61%  '; Text = $false }
    @{ Cell = 'D43'; Value = '12.49'; Text = $true }
    @{ Cell = 'E43'; Value = '  -6.83%  '; Text = $false }
    @{ Cell = 'D44'; Value = '15.18'; Text = $true }
    @{ Cell = 'E44'; Value = '  -9.18%  '; Text = $false }
    @{ Cell = 'D45'; Value = '0.6612'; Text = $true }
    @{ Cell = 'E45'; Value = '  -4.58%  '; Text = $false }
    @{ Cell = 'B46'; Value = 'Frax'; Text = $false }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'; Text = $false }
    @{ Cell = 'D46'; Value = '1.000'; Text = $true }
    @{ Cell = 'E46'; Value = '  +0.02%  '; Text = $false }
    @{ Cell = 'B47'; Value = 'NEARProtocol'; Text = $false }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; Text = $false }
    @{ Cell = 'D47'; Value = '2.313'; Text = $true }
    @{ Cell = 'E47'; Value = '  -5.10%  '; Text = $false }
    @{ Cell = 'D48'; Value = '3.990'; Text = $true }
    @{ Cell = 'E48'; Value = '  -2.60%  '; Text = $false }
    @{ Cell = 'D49'; Value = '132.53'; Text = $true }
    @{ Cell = 'E49'; Value = '  -1.39%  '; Text = $false }
    @{ Cell = 'D50'; Value = '0.07906'; Text = $true }
    @{ Cell = 'E50'; Value = '  -4.36%  '; Text = $false }
    @{ Cell = 'D51'; Value = '1.221'; Text = $true }
    @{ Cell = 'E51'; Value = '  -0.09%  '; Text = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Text) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
}
